$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 132, pushing the existing rows
# 132-187 down to 134-189 (the old last weekly block, previously rows
# 186-187, ends up as rows 188-189). The two freshly inserted rows
# (132-133) are then populated with a new weekly data block.
$ws.Range("A132:A133").EntireRow.Insert()

# Row 132 - "Primera" quality
$ws.Range("A132").Value = 11
$ws.Range("B132").Value = "Vega Monumental Concepción"
$ws.Range("C132").Value = "Bíobío"
$ws.Range("D132").Value = "2022-12-22"
$ws.Range("E132").Value = 8
$ws.Range("F132").Value = 100112044
$ws.Range("G132").Value = "Perejil"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 200
$ws.Range("K132").Value = 700
$ws.Range("L132").Value = 800
$ws.Range("M132").Value = 750
$ws.Range("N132").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O132").Value = "Región de Ñuble"
$ws.Range("P132").Value = 750
$ws.Range("Q132").Value = 1
$ws.Range("R132").Value = "Hortaliza"

# Row 133 - "Segunda" quality
$ws.Range("A133").Value = 11
$ws.Range("B133").Value = "Vega Monumental Concepción"
$ws.Range("C133").Value = "Bíobío"
$ws.Range("D133").Value = "2022-12-22"
$ws.Range("E133").Value = 8
$ws.Range("F133").Value = 100112044
$ws.Range("G133").Value = "Perejil"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Segunda"
$ws.Range("J133").Value = 100
$ws.Range("K133").Value = 600
$ws.Range("L133").Value = 600
$ws.Range("M133").Value = 600
$ws.Range("N133").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O133").Value = "Región de Ñuble"
$ws.Range("P133").Value = 600
$ws.Range("Q133").Value = 1
$ws.Range("R133").Value = "Hortaliza"
